$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D1').Value = 'Prequisites'
$ws.Range('E1').Value = ''
$ws.Range('D2').Value = 'applicants admission degree programme meet minimum general admission requirement either mean grade c+ kenya certificate secondary education ( kcse ) examination , 2 principal passes 1 subsidiary east african advanced certificate examination/ kenya advanced certificate education ( eace/ kace ) , equivalents.  grade c plain must diploma relevant field . addition , candidate must least c+ biology/biological sciences c+ agriculture kcse examinations . 3. holder diploma credit agricultural education related field college university recognized maseno university .'
$ws.Range('E2').Value = ''
$ws.Range('D4').Value = ' ii ) applicants must minimum c+ kcse equivalent .iv ) applicants may admitted strength diploma social sciences , arts humanities recognized institution .'
$ws.Range('E4').Value = ''
$ws.Range('D5').Value = ''
$ws.Range('D9').Value = 'c.s.e . minimum grade c mathematics two principal passes “ ” level/kenya advanced certificates education ( kace ) cue/ nec recognized diploma holders related areas eligible join first year exemptions . equivalent qualifications approved school business economics university senate .'
$ws.Range('E9').Value = ''
$ws.Range('D10').Value = ' obtained minimum grade c ( plain ) mathematics ksce equivalent .'
$ws.Range('E10').Value = ''
$ws.Range('D11').Value = ' candidate shall scored least ) grade c+ ( c plus ) english/ kiswahili . addition , candidate must attained least c ( c plain ) mathematics kcse equivalents examinations . candidate scores minimum principal pass chemistry , biology , geography subsidiary passes either maths physics a-level examinations eligible . candidates credit pass ordinary diploma hospitality tourism management related courses recognized institution qualify .'
$ws.Range('E11').Value = ''
$ws.Range('D12').Value = ' candidate shall scored least ) grade c+ ( c plus ) english/ kiswahili . addition , candidate must attained least c ( c plain ) mathematics kcse equivalents examinations . candidate scores minimum principal pass chemistry , biology , geography subsidiary passes either maths physics a-level examinations eligible . candidates credit pass ordinary diploma hospitality tourism management related courses recognized institution qualify .'
$ws.Range('E12').Value = ''
$ws.Range('D13').Value = 'ksce mean grade c+ , least c+ mathematics physics/physical science , kace least 2 principles passes mathematics physics 1 subsidiary pass , accredited diploma computer science /it/ related field least credit pass commission higher education accredited validated institution.  knec diploma holders computer studies least credit pass mean grade c kcse .'
$ws.Range('E13').Value = ''
$ws.Range('D15').Value = 'ksce mean grade c+ , least c+ mathematics physics/physical science , kace least 2 principles passes mathematics physics 1 subsidiary pass , accredited diploma computer science /it/ related field least credit pass commission higher education accredited validated institution.  knec diploma holders computer studies least credit pass mean grade c kcse .'
$ws.Range('E15').Value = ''
$ws.Range('D16').Value = 'minimum mean grade c+ kcse two principal passes kace/eaace diploma certificate ecde knec diploma education recognized institution'
$ws.Range('E16').Value = ''
